$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Generator Data")
$ws2 = $wb.Worksheets.Item("Yearly Fuel Costs")

# --- Sheet 1: "Generator Data" ---
# Update existing Nominal Capacity at upgrade 1 value
$ws1.Range("B6").Value = 114040.5301352051

# Shift the old "Investment at upgrade 1" / "Yearly O&M Cost at upgrade 1" /
# "Total actualized Fuel Cost" rows down, inserting the new "upgrade 2" / "upgrade 3"
# rows for Nominal Capacity, Investment, and Yearly O&M Cost.
$ws1.Range("A7").Value = "Nominal Capacity at upgrade 2"
$ws1.Range("B7").Value = 122757.35133520509

$ws1.Range("A8").Value = "Nominal Capacity at upgrade 3"
$ws1.Range("B8").Value = 131474.17253520511

$ws1.Range("A9").Value = "Investment at upgrade 1"
$ws1.Range("B9").Value = 47908.426709799664

$ws1.Range("A10").Value = "Investment at upgrade 2"
$ws1.Range("B10").Value = 3661.9365861200022

$ws1.Range("A11").Value = "Investment at upgrade 3"
$ws1.Range("B11").Value = 3661.9365861200022

$ws1.Range("A12").Value = "Yearly O&M Cost at upgrade 1"
$ws1.Range("B12").Value = 4790.8426709799669

$ws1.Range("A13").Value = "Yearly O&M Cost at upgrade 2"
$ws1.Range("B13").Value = 5157.0363295919669

$ws1.Range("A14").Value = "Yearly O&M Cost at upgrade 3"
$ws1.Range("B14").Value = 5523.229988203967

$ws1.Range("A15").Value = "Total actualized Fuel Cost"
$ws1.Range("B15").Value = 740292.51719902863

# Apply same style (s="1", bold/boxed header label) as the rest of column A
# to the newly-written cells
$ws1.Range("A9").Copy()
$ws1.Range("A10:A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Sheet 2: "Yearly Fuel Costs" ---
# Overwrite the first three data rows with new (Salvage Value / Battery
# Replacement driven) figures, then delete the now-unused rows 5:21.
$ws2.Range("B2").Value = 220530.57755536071
$ws2.Range("B3").Value = 246612.91507201109
$ws2.Range("B4").Value = 273149.02457167441

$ws2.Rows("5:21").Delete()
